$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Title (Heading1): drop trailing " dos Materiais"
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "LOM3256 -  Tópicos em Cálculo de Estrutura Eletrônica dos Materiais",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LOM3256 -  Tópicos em Cálculo de Estrutura Eletrônica", 2) | Out-Null

# ---------------------------------------------------------------
# 2. Subtitle (Heading3): drop trailing " of materials"
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Methods of electronic structure calculation of materials",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Methods of electronic structure calculation", 2) | Out-Null

# ---------------------------------------------------------------
# 3. Ativação date update
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Ativação: 15/07/2015",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ativação: 01/01/2023", 2) | Out-Null

# ---------------------------------------------------------------
# 4. "Objetivos" paragraph: collapse the manual line breaks into a
#    single run, then add a new italic English translation
#    paragraph right after it.
# ---------------------------------------------------------------
$rngObjetivos = $d.Content
$rngObjetivos.Find.Execute("Objetivos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pObjetivosHeading = $rngObjetivos.Paragraphs(1)
$pObjetivos = $pObjetivosHeading.Next()
$pObjetivos.Range.Text = "Propiciar ao aluno uma visão básica sobre os principais métodos de determinação teórica da estrutura eletrônica dos materiais, com enfoque em sólidos cristalinos, mas também em materiais bidimensionais e nanoestruturados.O principal método de cálculo a ser empregado no curso será a Teoria do Funcional da Densidade(Density Functional Theory, DFT), em algumas de suas muitas variantes. Ao final do curso, o aluno estará apto a determinar propriedades dos materiais como estruturas de bandas, densidades de estados, superfícies de Fermi e constantes elásticas, usando um ou mais dos métodos e códigos computacionais apresentados em aula."

$pObjetivos.Range.InsertParagraphAfter()
$pObjetivosEn = $pObjetivos.Next()
$objetivosEnText = "Provide the student with a basic view of the main methods of theoretical determination of the electronic structure, focusing on crystalline solids, but also on molecules, two-dimensional materials and nanostructured materials. The main calculation method to be used in the course will be the Density Functional Theory (DFT), in some of its many variants. At the end of the course, the student will be able to determine material properties such as band structures, densities of states, elastic constants, and Fermi surfaces, using one or more of the methods and computer codes presented in class."
$startObjEn = $pObjetivosEn.Range.Start
$pObjetivosEn.Range.Text = $objetivosEnText
$d.Range($startObjEn, $startObjEn + $objetivosEnText.Length).Italic = 1

# ---------------------------------------------------------------
# 5. "Programa resumido" paragraph: collapse breaks, add English
#    italic translation paragraph right after it.
# ---------------------------------------------------------------
$rngResumido = $d.Content
$rngResumido.Find.Execute("Programa resumido", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pResumidoHeading = $rngResumido.Paragraphs(1)
$pResumido = $pResumidoHeading.Next()
$pResumido.Range.Text = "Revisão de mecânica quântica; Revisão de física do estado sólido; Método de Hartree-Fock; Teoria do funcional da densidade; Métodos de ondas planas e pseudo-potenciais; Códigos computacionais"

$pResumido.Range.InsertParagraphAfter()
$pResumidoEn = $pResumido.Next()
$resumidoEnText = "Review of Quantum Mechanics; Review of Solid State Physics; Hartree-Fock Method; Density Functional Theory; Plane and pseudopotential wave methods; computer codes"
$startResEn = $pResumidoEn.Range.Start
$pResumidoEn.Range.Text = $resumidoEnText
$d.Range($startResEn, $startResEn + $resumidoEnText.Length).Italic = 1

# ---------------------------------------------------------------
# 6. "Programa" detailed paragraph (the heading that immediately
#    follows the "Programa resumido" body -- search begins right
#    after the Resumido section to avoid matching "Programa
#    resumido" itself). Collapse breaks (no spaces inserted at
#    all), add English italic translation paragraph right after.
# ---------------------------------------------------------------
$rngPrograma = $d.Range($pResumidoEn.Range.End, $d.Content.End)
$rngPrograma.Find.Execute("Programa", $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pProgramaHeading = $rngPrograma.Paragraphs(1)
$pPrograma = $pProgramaHeading.Next()
$pPrograma.Range.Text = "Revisão de mecânica quânticao Equação de Schrödingero Átomo do hidrogênio e orbitais atômicoso Notação de Diraco Princípio variacionalo Combinação linear de orbitais atômicosRevisão de física do estado sólidoo Espaço direto e recíprocoo Teorema de Blocho Zona de Brillouino Bandas de energia e densidade de estadoso Energia de Fermi e superfície de Fermio Aproximação de elétrons livresMétodo de Hartree-Focko Determinantes de Slatero Equação de Hartree-Focko Potencial de troca e correlaçãoo Algoritmo autoconsistenteTeoria do funcional da densidadeo Teoremas de Hohenberg-Kohno Equações de Kohn-Shamo Funcionais de troca e correlação: LDA, GGA, etc.Métodos de ondas planas e pseudo-potenciaiso Bases de ondas planaso Pseudo-potenciaiso Bases de ondas planas aumentadas e linearizadaso Método FP-LAPWCódigos computacionaiso Quantum Espressoo Elko Wien2ko VASP"

$pPrograma.Range.InsertParagraphAfter()
$pProgramaEn = $pPrograma.Next()
$programaEnText = "• Review of quantum mechanics: Schrödinger's equation; Hydrogen atom and atomic orbitals; Dirac notation; Variational principle; Linear combination of atomic orbitals. • Solid state physics review: Direct and reciprocal space; Bloch's Theorem; Brillouin zone; Energy bands and density of states; Fermi energy and Fermi surface; Free electrons Approximation. • Hartree-Fock method: Slater determinants; Hartree-Fock equation; Exchange and correlation potential; Self-consistent algorithm. • Density functional theory: Hohenberg-Kohn theorems; Kohn-Sham equations; Exchange and correlation functionals: LDA, GGA, etc. • Plane and pseudopotential wave methods: Plane wave bases; Pseudo-potentials; • Augmented and linearized plane wave bases: FP-LAPW method. • Computer codes: NWCHEM, Quantum Espresso, , Wien2k, exciting, VASP, etc."
$startProgEn = $pProgramaEn.Range.Start
$pProgramaEn.Range.Text = $programaEnText
$d.Range($startProgEn, $startProgEn + $programaEnText.Length).Italic = 1

# ---------------------------------------------------------------
# 7. Requisitos: swap course requirement
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "LOM3215 -  Física do Estado Sólido  (Requisito)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "LOM3226 -  Mecânica Quântica  (Requisito)", 2) | Out-Null

Write-Host "All edits applied"
